$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '37.118.38'
Set-TextValue 'E2' '  -1.42%  '

# Row 3
Set-TextValue 'D3' '2.090.30'
Set-TextValue 'E3' '  +7.60%  '

# Row 4
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.14%  '

# Row 5
Set-TextValue 'D5' '249.98'
Set-TextValue 'E5' '  -0.77%  '

# Row 6
Set-TextValue 'D6' '0.651'
Set-TextValue 'E6' '  -6.15%  '

# Row 7
Set-TextValue 'E7' '  +0.21%  '

# Row 8
Set-TextValue 'D8' '50.00'
Set-TextValue 'E8' '  +3.09%  '

# Row 9
Set-TextValue 'D9' '60.64'
Set-TextValue 'E9' '  +2.94%  '

# Row 10
Set-TextValue 'D10' '0.367'
Set-TextValue 'E10' '  -4.11%  '

# Row 11
Set-TextValue 'D11' '0.0738'
Set-TextValue 'E11' '  -4.59%  '

# Row 12
Set-TextValue 'E12' '  +4.16%  '

# Row 13
Set-TextValue 'D13' '15.15'
Set-TextValue 'E13' '  -4.35%  '

# Row 14
Set-TextValue 'D14' '2.375.25'
Set-TextValue 'E14' '  +6.75%  '

# Row 15
Set-TextValue 'D15' '0.828'
Set-TextValue 'E15' '  -2.32%  '

# Row 16
Set-TextValue 'D16' '2.082.92'
Set-TextValue 'E16' '  +7.26%  '

# Row 17
Set-TextValue 'D17' '5.05'
Set-TextValue 'E17' '  -2.91%  '

# Row 18
Set-TextValue 'D18' '37.022.44'
Set-TextValue 'E18' '  -1.75%  '

# Row 19
Set-TextValue 'D19' '71.88'
Set-TextValue 'E19' '  -5.06%  '

# Row 20
Set-TextValue 'D20' '0.0₃0818'
Set-TextValue 'E20' '  -5.56%  '

# Row 21
Set-TextValue 'D21' '13.18'
Set-TextValue 'E21' '  -4.14%  '

# Row 22
Set-TextValue 'D22' '238.26'
Set-TextValue 'E22' '  -6.38%  '

# Row 23
Set-TextValue 'D23' '5.17'
Set-TextValue 'E23' '  -1.45%  '

# Row 24
Set-TextValue 'E24' '  +0.10%  '

# Row 25
Set-TextValue 'E25' '  -2.32%  '

# Row 26
Set-TextValue 'D26' '168.58'
Set-TextValue 'E26' '  -0.58%  '

# Row 27
Set-TextValue 'D27' '9.24'
Set-TextValue 'E27' '  +2.76%  '

# Row 28
Set-TextValue 'D28' '20.73'
Set-TextValue 'E28' '  +9.15%  '

# Row 29
Set-TextValue 'D29' '1.99'
Set-TextValue 'E29' '  -6.46%  '

# Row 30
Set-TextValue 'E30' '  -6.13%  '

# Row 31
Set-TextValue 'D31' '22.52'
Set-TextValue 'E31' '  +14.33%  '

# Row 32
Set-TextValue 'D32' '1.07'
Set-TextValue 'E32' '  +18.67%  '

# Row 33
Set-TextValue 'D33' '4.46'
Set-TextValue 'E33' '  -3.71%  '

# Row 34
Set-TextValue 'D34' '0.0602'
Set-TextValue 'E34' '  -2.68%  '

# Row 35
Set-TextValue 'E35' '  -2.50%  '

# Row 36
Set-TextValue 'E36' '  +0.00%  '

# Row 37
Set-TextValue 'B37' 'LidoDAOToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D37' '2.26'
Set-TextValue 'E37' '  +12.48%  '

# Row 38
Set-TextValue 'B38' 'WEMIXToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D38' '1.82'
Set-TextValue 'E38' '  -3.26%  '

# Row 39
Set-TextValue 'D39' '4.04'

# Row 40
Set-TextValue 'E40' '  -10.22%  '

# Row 41
Set-TextValue 'D41' '17.49'
Set-TextValue 'E41' '  -0.84%  '

# Row 42
Set-TextValue 'D42' '0.0221'
Set-TextValue 'E42' '  -3.49%  '

# Row 43
Set-TextValue 'E43' '  +1.73%  '

# Row 44
Set-TextValue 'D44' '97.57'
Set-TextValue 'E44' '  -7.71%  '

# Row 45
Set-TextValue 'D45' '2.77'
Set-TextValue 'E45' '  -4.02%  '

# Row 46
Set-TextValue 'D46' '0.0877'
Set-TextValue 'E46' '  +3.67%  '

# Row 47
Set-TextValue 'D47' '2.96'
Set-TextValue 'E47' '  +4.90%  '

# Row 48
Set-TextValue 'D48' '1.302.13'
Set-TextValue 'E48' '  -4.23%  '

# Row 49
Set-TextValue 'D49' '6.83'
Set-TextValue 'E49' '  +5.64%  '

# Row 50
Set-TextValue 'D50' '2.254.78'
Set-TextValue 'E50' '  +6.61%  '

# Row 51
Set-TextValue 'D51' '2.24'
Set-TextValue 'E51' '  -8.29%  '
